{"js": "// Replace each three-digit-division-by-one-digit expression in the\n// document's table cells with its new value, as described by the diff.\n// Every <w:t> run that holds an equation like \"160\u00f74=\" changes to a new\n// equation (\"240\u00f74=\" etc.) \u2014 the date heading and blank cells are left\n// untouched. Because every \"before\" value is unique in the document, a\n// simple search + replace per pair is safe and unambiguous.\n\nconst replacements = [\n  [\"160\u00f74=\", \"240\u00f74=\"],\n  [\"642\u00f74=\", \"317\u00f72=\"],\n  [\"657\u00f77=\", \"224\u00f75=\"],\n  [\"569\u00f77=\", \"450\u00f72=\"],\n  [\"155\u00f76=\", \"821\u00f76=\"],\n  [\"778\u00f79=\", \"112\u00f72=\"],\n  [\"654\u00f76=\", \"711\u00f72=\"],\n  [\"858\u00f75=\", \"254\u00f75=\"],\n  [\"904\u00f72=\", \"469\u00f76=\"],\n  [\"743\u00f76=\", \"570\u00f77=\"],\n  [\"777\u00f79=\", \"341\u00f72=\"],\n  [\"580\u00f74=\", \"618\u00f75=\"],\n  [\"915\u00f72=\", \"341\u00f78=\"],\n  [\"662\u00f74=\", \"558\u00f75=\"],\n  [\"483\u00f73=\", \"997\u00f72=\"],\n  [\"443\u00f72=\", \"815\u00f72=\"],\n  [\"398\u00f76=\", \"959\u00f79=\"],\n  [\"261\u00f77=\", \"995\u00f76=\"],\n  [\"606\u00f74=\", \"153\u00f75=\"],\n  [\"285\u00f78=\", \"343\u00f73=\"],\n  [\"310\u00f76=\", \"589\u00f78=\"],\n  [\"491\u00f77=\", \"737\u00f78=\"],\n  [\"922\u00f76=\", \"856\u00f72=\"],\n  [\"529\u00f74=\", \"252\u00f77=\"],\n  [\"469\u00f72=\", \"352\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-division-by-one-digit expression in the\n# document's table cells with its new value, as described by the diff.\n# Every text run that holds an equation like \"160\u00f74=\" changes to a new\n# equation (\"240\u00f74=\" etc.) \u2014 the date heading and blank cells are left\n# untouched. Because every \"before\" value is unique in the document, a\n# simple Find/Replace per pair (restricted to an exact, case-sensitive,\n# non-wildcard match) is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"160\u00f74=\", \"240\u00f74=\"),\n  @(\"642\u00f74=\", \"317\u00f72=\"),\n  @(\"657\u00f77=\", \"224\u00f75=\"),\n  @(\"569\u00f77=\", \"450\u00f72=\"),\n  @(\"155\u00f76=\", \"821\u00f76=\"),\n  @(\"778\u00f79=\", \"112\u00f72=\"),\n  @(\"654\u00f76=\", \"711\u00f72=\"),\n  @(\"858\u00f75=\", \"254\u00f75=\"),\n  @(\"904\u00f72=\", \"469\u00f76=\"),\n  @(\"743\u00f76=\", \"570\u00f77=\"),\n  @(\"777\u00f79=\", \"341\u00f72=\"),\n  @(\"580\u00f74=\", \"618\u00f75=\"),\n  @(\"915\u00f72=\", \"341\u00f78=\"),\n  @(\"662\u00f74=\", \"558\u00f75=\"),\n  @(\"483\u00f73=\", \"997\u00f72=\"),\n  @(\"443\u00f72=\", \"815\u00f72=\"),\n  @(\"398\u00f76=\", \"959\u00f79=\"),\n  @(\"261\u00f77=\", \"995\u00f76=\"),\n  @(\"606\u00f74=\", \"153\u00f75=\"),\n  @(\"285\u00f78=\", \"343\u00f73=\"),\n  @(\"310\u00f76=\", \"589\u00f78=\"),\n  @(\"491\u00f77=\", \"737\u00f78=\"),\n  @(\"922\u00f76=\", \"856\u00f72=\"),\n  @(\"529\u00f74=\", \"252\u00f77=\"),\n  @(\"469\u00f72=\", \"352\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
